$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "PRIMERA PREGUNTA" header to "PRIMERA PREGUNTA1"
$ws.Range("F1").Value = "PRIMERA PREGUNTA1"

# Add the two new result columns used by the robot run: ESTADO / MENSAJE
$ws.Range("K1").Value = "ESTADO"
$ws.Range("L1").Value = "MENSAJE"

# Give column L the wider width needed to show a full message
$ws.Columns.Item(12).ColumnWidth = 43

# Clear the stray duplicated answer in H3 (second row no longer repeats RESPUESTA 1)
$ws.Range("H3").ClearContents()

# The new K/L output cells for both data rows reuse the plain (no-underline)
# formatting that used to live only on F3 - copy it across, then strip the
# now-redundant underline formatting back off of F3 itself.
$ws.Range("F3").Copy() | Out-Null
$ws.Range("K2:L3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").ClearFormats()

# Leave the selection where the robot left it after writing the new columns
$ws.Range("K2:L3").Select() | Out-Null
